$d = $word.ActiveDocument

# --- Paragraph: "Utilicen el léxico..." ---
# "...para que, además de expresiones aritméticas (que ya acepta), sea capaz de
#  reconocer asignaciones y expresiones condicionales." ->
# "...para que acepte instrucciones más complejas que las expresiones que ya acepta."
$old2 = "para que, además de expresiones aritméticas (que ya acepta), sea capaz de reconocer asignaciones y expresiones condicionales."
$new2 = "para que acepte instrucciones más complejas que las expresiones que ya acepta."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- Paragraph: "Toda instrucción válida..." ---
$old3 = "Las asignaciones se forman con un identificador seguido de un operador de asignación (=) seguido por una expresión matemática. Las expresiones matemáticas pueden ser expresiones aritméticas o expresiones condicionales, escritas entre paréntesis."
$new3 = "Las asignaciones se forman con un identificador seguido de un operador de asignación (=) seguido por una expresión. Las expresiones pueden ser operaciones aritméticas binarias entre expresiones, expresiones entre paréntesis, expresiones condicionales, constantes, o identificadores."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# --- Paragraph: "Una expresión condicional..." ---
$old4a = "Una expresión condicional incluye"
$new4a = "Una expresión condicional se escribe entre llaves {} e incluye"
$d.Content.Find.Execute($old4a, $true, $false, $false, $false, $false, $true, 1, $false, $new4a, 2) | Out-Null

$old4b = "separadas con el símbolo :"
$new4b = "separadas con el símbolo :."
$d.Content.Find.Execute($old4b, $true, $false, $false, $false, $false, $true, 1, $false, $new4b, 2) | Out-Null

# --- Example expressions: remove / restyle parentheses into brace ({}) conditional syntax ---
$d.Content.Find.Execute("c = (a + e) `$", $true, $false, $false, $false, $false, $true, 1, $false, "c = a + e `$", 2) | Out-Null

$d.Content.Find.Execute("v = (5.3 * (2 + x)) `$", $true, $false, $false, $false, $false, $true, 1, $false, "v = 5.3 * (2 + x) `$", 2) | Out-Null

$d.Content.Find.Execute("y = ((x + 1) > (z) ? (2.34) : (z – 4)) `$", $true, $false, $false, $false, $false, $true, 1, $false, "y = {x + 1 > z ? 2.34 : z – 4} `$", 2) | Out-Null

$d.Content.Find.Execute("tasa = (4 * ((tasa) <= (0.2) ? (5) : (8)) `$", $true, $false, $false, $false, $false, $true, 1, $false, "tasa = 4 * {tasa <= 0.2 ? 5 : 8} `$", 2) | Out-Null

$d.Content.Find.Execute("x = (((x) != (y/3) ? (x+1) : (y-2)) == (z) ? (2.2) : (3.3)) `$", $true, $false, $false, $false, $false, $true, 1, $false, "x = {{x != y/3 ? (x+1) * 1.2 : y-2} == z ? 2.2 : {x >= y ? 2 : 4}} `$", 2) | Out-Null

# --- Indent the 5 example paragraphs by 708 twips (35.4 pt) ---
for ($i = 7; $i -le 11; $i++) {
    $p = $d.Paragraphs($i)
    $p.LeftIndent = 35.4
}
